$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "A: " -> "A: Krishna"   (unique occurrence; ReplaceOne avoids
#    touching the other "A: " runs later in the document)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "A: ", $true, $false, $false, $false, $false, $true, 1, $false,
    "A: Krishna", 1) | Out-Null

# ---------------------------------------------------------------------
# 2) "B:" / "C:" / "D:" run (with line breaks) -> "B: Ihsan" / "C: Pranav " / "D:avid Kim"
#    Scope the Find to the containing paragraph so we don't touch the
#    later "B: Or, had ..." / "C: These problems ..." / "D: So we asked ..." paragraphs.
# ---------------------------------------------------------------------
$pBCD = $d.Paragraphs.Item(28)
$pBCD.Range.Find.Execute(
    "B:", $true, $false, $false, $false, $false, $true, 1, $false,
    "B: Ihsan", 1) | Out-Null

$pBCD = $d.Paragraphs.Item(28)
$pBCD.Range.Find.Execute(
    "C:", $true, $false, $false, $false, $false, $true, 1, $false,
    "C: Pranav ", 1) | Out-Null

$pBCD = $d.Paragraphs.Item(28)
$pBCD.Range.Find.Execute(
    "D:", $true, $false, $false, $false, $false, $true, 1, $false,
    "D:avid Kim", 1) | Out-Null

# ---------------------------------------------------------------------
# 3) The "A: When was ... / B: Or, had ... / C: These problems ... /
#    D: So we asked ..." block gets a blank paragraph inserted before
#    each of B/C/D, and the A and B lines get reworded.
# ---------------------------------------------------------------------

# 3a. Reword the "A: When was the last time ..." line.
$d.Content.Find.Execute(
    "A: When was the last time you forgot, or too lazy to turn off your lights ? ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A: When was the last time you forgot to turn off your lights? ", 1) | Out-Null

# 3b. Insert a blank paragraph right before "B: Or, had ..." then reword it.
$rngB = $d.Content
$rngB.Find.Execute(
    "B: Or, had to get to class really quickly, leaving no time to turn off all, or any of your electronics? ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngB.InsertParagraphBefore()

$d.Content.Find.Execute(
    "B: Or, had to get to class really quickly, leaving no time to turn off all, or any of your electronics? ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "B: Or had to get to class quickly, leaving no time to turn off your electronics? ", 1) | Out-Null

# 3c. Insert a blank paragraph right before "C: These problems ..." (text unchanged).
$rngC = $d.Content
$rngC.Find.Execute(
    "C: These problems are common. But the products to help you are not.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngC.InsertParagraphBefore()

# 3d. Insert a blank paragraph right before "D: So we asked ..." (text unchanged).
$rngD = $d.Content
$rngD.Find.Execute(
    "D: So we asked, why not develop something amazing? Why not make something that is green? Why not find out where you use up most of your energy? So, we created OutLite, a smartplug that incorporates all of this.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngD.InsertParagraphBefore()

# ---------------------------------------------------------------------
# 4) The trailing "b: " paragraph becomes blank, and three new
#    speaker paragraphs (each preceded by a blank paragraph) follow it.
#    Locate the anchor paragraph ("A: We have designed OutLite ...")
#    by search, then address everything after it by index so we don't
#    depend on hard-coded absolute paragraph numbers.
# ---------------------------------------------------------------------
$rngAnchor = $d.Content
$rngAnchor.Find.Execute(
    "A: We have designed OutLite", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$idxDesigned = $rngAnchor.Paragraphs.Item(1).Index
$idxB = $idxDesigned + 1

# Clear the original "b: " text, leaving the paragraph blank.
$d.Paragraphs.Item($idxB).Range.Text = ""

# Insert "B: By completely ..." right after the (now blank) paragraph.
$d.Paragraphs.Item($idxB).Range.InsertParagraphAfter()
$idxBText = $idxB + 1
$d.Paragraphs.Item($idxBText).Range.Text = "B: By completely cutting off the electricity, you can additionally save on standby power usage, which approximately comprises 10% of your electricity bill."

# Blank paragraph, then "C: So far ..."
$d.Paragraphs.Item($idxBText).Range.InsertParagraphAfter()
$idxBlank2 = $idxBText + 1

$d.Paragraphs.Item($idxBlank2).Range.InsertParagraphAfter()
$idxCText = $idxBlank2 + 1
$d.Paragraphs.Item($idxCText).Range.Text = "C: So far, we have launched a beta program, where we are inviting potential customers to join us in testing the OutLite. You can sign up for the beta on facebook. We plan to call on our beta-testers to figure out how the smartplug will be used specifically, and to fine-tune the features currently in development."

# Blank paragraph, then "D: We are currently ..."
$d.Paragraphs.Item($idxCText).Range.InsertParagraphAfter()
$idxBlank3 = $idxCText + 1

$d.Paragraphs.Item($idxBlank3).Range.InsertParagraphAfter()
$idxDText = $idxBlank3 + 1
$d.Paragraphs.Item($idxDText).Range.Text = "D: We are currently researching and developing the prototype using the SparkCore, a wi-fi development platform. After we successfully develop the prototype, we will design our smartplug's enclosure with SolidWorks, incorporate bluetooth technology for a lower price and a smaller carbon footprint. Thank you for listening, and turn down for Watt!"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
